# Apply scheduled-runner updates to profit/price columns across sheets.
# Values below come from the latest Universalis price pull; only the
# numeric columns (H-N) change, row/item identity columns are untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 626.7692
$ws.Range("J92").Value = 5000
$ws.Range("L92").Value = 5000
$ws.Range("N92").Value = -7496

$ws.Range("H137").Value = 1235.3684
$ws.Range("I137").Value = 718
$ws.Range("J137").Value = 1332.375
$ws.Range("K137").Value = 2154
$ws.Range("L137").Value = 3997.125
$ws.Range("M137").Value = 396
$ws.Range("N137").Value = -9097.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2480.3635
$ws.Range("I2").Value = 2401
$ws.Range("J2").Value = 2692
$ws.Range("K2").Value = 2401
$ws.Range("L2").Value = 2692
$ws.Range("M2").Value = -2288
$ws.Range("N2").Value = -2918

$ws.Range("H110").Value = 1060.6666
$ws.Range("I110").Value = 752.4138
$ws.Range("J110").Value = 10000
$ws.Range("K110").Value = 752.4138
$ws.Range("L110").Value = 10000
$ws.Range("M110").Value = 1292.5862
$ws.Range("N110").Value = -14090

$ws.Range("H116").Value = 2480.3635
$ws.Range("I116").Value = 2401
$ws.Range("J116").Value = 2692
$ws.Range("K116").Value = 2401
$ws.Range("L116").Value = 2692
$ws.Range("M116").Value = -107
$ws.Range("N116").Value = -7280

$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -59178

$ws.Range("H119").Value = 46666.668
$ws.Range("J119").Value = 46666.668
$ws.Range("L119").Value = 46666.668
$ws.Range("N119").Value = -56342.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2480.3635
$ws.Range("I3").Value = 2401
$ws.Range("J3").Value = 2692
$ws.Range("K3").Value = 2401
$ws.Range("L3").Value = 2692
$ws.Range("M3").Value = -2287
$ws.Range("N3").Value = -2920

$ws.Range("H20").Value = 4017.5557
$ws.Range("I20").Value = 3755.077
$ws.Range("J20").Value = 4700
$ws.Range("K20").Value = 3755.077
$ws.Range("L20").Value = 4700
$ws.Range("M20").Value = -3508.077
$ws.Range("N20").Value = -5194

$ws.Range("H63").Value = 27635.5
$ws.Range("J63").Value = 27635.5
$ws.Range("L63").Value = 27635.5
$ws.Range("N63").Value = -29007.5

$ws.Range("H66").Value = 27635.5
$ws.Range("J66").Value = 27635.5
$ws.Range("L66").Value = 82906.5
$ws.Range("N66").Value = -89770.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 901
$ws.Range("I132").Value = 1002
$ws.Range("J132").Value = 800
$ws.Range("K132").Value = 9018
$ws.Range("L132").Value = 7200
$ws.Range("M132").Value = -6488
$ws.Range("N132").Value = -12260

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2400.1785
$ws.Range("I80").Value = 2216.6667
$ws.Range("J80").Value = 2537.8125
$ws.Range("K80").Value = 2216.6667
$ws.Range("L80").Value = 2537.8125
$ws.Range("M80").Value = -1218.6667
$ws.Range("N80").Value = -4533.8125

$ws.Range("H83").Value = 2400.1785
$ws.Range("I83").Value = 2216.6667
$ws.Range("J83").Value = 2537.8125
$ws.Range("K83").Value = 11083.3335
$ws.Range("L83").Value = 12689.0625
$ws.Range("M83").Value = -6091.333500000001
$ws.Range("N83").Value = -22673.0625

$ws.Range("H92").Value = 10000
$ws.Range("J92").Value = 10000
$ws.Range("L92").Value = 10000
$ws.Range("N92").Value = -13744

$ws.Range("H93").Value = 27559.25
$ws.Range("J93").Value = 27559.25
$ws.Range("L93").Value = 27559.25
$ws.Range("N93").Value = -31303.25

$ws.Range("H95").Value = 53500
$ws.Range("J95").Value = 53500
$ws.Range("L95").Value = 53500
$ws.Range("N95").Value = -58992

$ws.Range("H96").Value = 19490.25
$ws.Range("J96").Value = 19490.25
$ws.Range("L96").Value = 19490.25
$ws.Range("N96").Value = -24982.25

$ws.Range("H98").Value = 30321.5
$ws.Range("J98").Value = 30321.5
$ws.Range("L98").Value = 30321.5
$ws.Range("N98").Value = -36311.5

$ws.Range("H99").Value = 13703
$ws.Range("I99").Value = 12986.833
$ws.Range("J99").Value = 18000
$ws.Range("K99").Value = 12986.833
$ws.Range("L99").Value = 18000
$ws.Range("M99").Value = -10740.833
$ws.Range("N99").Value = -22492

$ws.Range("H101").Value = 61000
$ws.Range("J101").Value = 61000
$ws.Range("L101").Value = 61000
$ws.Range("N101").Value = -67490

$ws.Range("H102").Value = 3001.2964
$ws.Range("I102").Value = 2224.1667
$ws.Range("J102").Value = 4555.5557
$ws.Range("K102").Value = 2224.1667
$ws.Range("L102").Value = 4555.5557
$ws.Range("M102").Value = -602.1667000000002
$ws.Range("N102").Value = -7799.5557

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2060.7576
$ws.Range("I136").Value = 1371.9524
$ws.Range("J136").Value = 3266.1667
$ws.Range("K136").Value = 4115.857199999999
$ws.Range("L136").Value = 9798.500100000001
$ws.Range("M136").Value = -1565.857199999999
$ws.Range("N136").Value = -14898.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 123433.664
$ws.Range("J80").Value = 123433.664
$ws.Range("L80").Value = 123433.664
$ws.Range("N80").Value = -125429.664

$ws.Range("H83").Value = 123433.664
$ws.Range("J83").Value = 123433.664
$ws.Range("L83").Value = 370300.992
$ws.Range("N83").Value = -380284.992

$ws.Range("H92").Value = 39850
$ws.Range("J92").Value = 39850
$ws.Range("L92").Value = 39850
$ws.Range("N92").Value = -44842

$ws.Range("H94").Value = 34833.332
$ws.Range("J94").Value = 34833.332
$ws.Range("L94").Value = 34833.332
$ws.Range("N94").Value = -36635.332

$ws.Range("H95").Value = 36355.145
$ws.Range("I95").Value = 41000
$ws.Range("J95").Value = 35581
$ws.Range("K95").Value = 41000
$ws.Range("L95").Value = 35581
$ws.Range("M95").Value = -38254
$ws.Range("N95").Value = -41073

$ws.Range("H96").Value = 2234.3333
$ws.Range("I96").Value = 2234.3333
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 2234.3333
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -861.3332999999998
$ws.Range("N96").ClearContents()

$ws.Range("H98").Value = 35750
$ws.Range("J98").Value = 35750
$ws.Range("L98").Value = 35750
$ws.Range("N98").Value = -41740

$ws.Range("H100").Value = 5682590.5
$ws.Range("I100").Value = 6494282
$ws.Range("J100").Value = 750
$ws.Range("K100").Value = 12988564
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = -12988023
$ws.Range("N100").Value = -2582

$ws.Range("H104").Value = 22090
$ws.Range("J104").Value = 22090
$ws.Range("L104").Value = 22090
$ws.Range("N104").Value = -29078

$ws.Range("H110").Value = 39996.668
$ws.Range("J110").Value = 39996.668
$ws.Range("L110").Value = 39996.668
$ws.Range("N110").Value = -48176.668

$ws.Range("H111").Value = 46383.75
$ws.Range("J111").Value = 46383.75
$ws.Range("L111").Value = 46383.75
$ws.Range("N111").Value = -54563.75

$ws.Range("H113").Value = 425.25
$ws.Range("I113").Value = 267
$ws.Range("J113").Value = 900
$ws.Range("K113").Value = 801
$ws.Range("L113").Value = 2700
$ws.Range("M113").Value = 1369
$ws.Range("N113").Value = -7040

$ws.Range("H114").Value = 39333.332
$ws.Range("J114").Value = 39333.332
$ws.Range("L114").Value = 39333.332
$ws.Range("N114").Value = -48011.332

$ws.Range("H115").Value = 50000
$ws.Range("J115").Value = 50000
$ws.Range("L115").Value = 50000
$ws.Range("N115").Value = -53134

$ws.Range("H118").Value = 46748.75
$ws.Range("J118").Value = 46748.75
$ws.Range("L118").Value = 46748.75
$ws.Range("N118").Value = -50062.75

$ws.Range("H119").Value = 50698
$ws.Range("J119").Value = 50698
$ws.Range("L119").Value = 50698
$ws.Range("N119").Value = -60374

$ws.Range("H121").Value = 42236.668
$ws.Range("J121").Value = 42236.668
$ws.Range("L121").Value = 42236.668
$ws.Range("N121").Value = -45730.668

$ws.Range("H132").Value = 1213.4062
$ws.Range("I132").Value = 938.4737
$ws.Range("J132").Value = 1615.2307
$ws.Range("K132").Value = 2815.4211
$ws.Range("L132").Value = 4845.6921
$ws.Range("M132").Value = -285.4211
$ws.Range("N132").Value = -9905.6921

$ws.Range("H136").Value = 3470.5957
$ws.Range("I136").Value = 1216.4166
$ws.Range("J136").Value = 5822.7827
$ws.Range("K136").Value = 3649.2498
$ws.Range("L136").Value = 17468.3481
$ws.Range("M136").Value = -1099.2498
$ws.Range("N136").Value = -22568.3481

